$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 23, pushing existing rows 23..122 down to 24..123.
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the new weekly data entry.
$ws.Cells.Item(23, 1).Value = 3
$ws.Cells.Item(23, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value = [DateTime]"2021-12-28"
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 5
$ws.Cells.Item(23, 6).Value = 100112026
$ws.Cells.Item(23, 7).Value = "Haba"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 115
$ws.Cells.Item(23, 11).Value = 7500
$ws.Cells.Item(23, 12).Value = 8000
$ws.Cells.Item(23, 13).Value = 7739
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(23, 16).Value = 310
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
